# Applies the cryptos list refresh described in the commit:
# "Updated cryptos list on Sun Apr 21 13:36:23 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '65.105.97'
    'E2' = '  +1.84%  '
    'D3' = '3.164.21'
    'E3' = '  +3.25%  '
    'E4' = '  +0.09%  '
    'D5' = '577.28'
    'E5' = '  +3.25%  '
    'D6' = '150.46'
    'E6' = '  +5.27%  '
    'D8' = '3.161.69'
    'E8' = '  +3.26%  '
    'D9' = '0.528'
    'E9' = '  +2.12%  '
    'E10' = '  +4.34%  '
    'E11' = '  -0.56%  '
    'D12' = '0.501'
    'E12' = '  +4.25%  '
    'E13' = '  +13.89%  '
    'D14' = '37.28'
    'E14' = '  +5.57%  '
    'D15' = '3.682.85'
    'E15' = '  +3.35%  '
    'D16' = '65.212.66'
    'E16' = '  +2.00%  '
    'B17' = 'WrappedEther'
    'C17' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    'D17' = '3.164.51'
    'E17' = '  +3.33%  '
    'B18' = 'Polkadot'
    'C18' = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
    'D18' = '7.14'
    'E18' = '  +5.16%  '
    'E19' = '  +1.21%  '
    'D20' = '510.51'
    'E20' = '  +4.38%  '
    'D21' = '14.86'
    'E21' = '  +4.10%  '
    'D22' = '15.38'
    'E22' = '  +4.76%  '
    'E23' = '  +4.86%  '
    'E24' = '  +3.16%  '
    'D25' = '84.66'
    'E25' = '  +2.14%  '
    'E26' = '  -0.05%  '
    'D27' = '2.92'
    'E27' = '  +4.01%  '
    'D28' = '8.91'
    'E28' = '  +9.50%  '
    'E29' = '  +6.50%  '
    'D30' = '27.81'
    'E30' = '  +4.83%  '
    'D31' = '2.77'
    'E31' = '  +9.50%  '
    'D32' = '1.00'
    'E32' = '  +0.15%  '
    'D33' = '1.19'
    'E33' = '  +2.90%  '
    'D34' = '6.27'
    'E34' = '  +9.69%  '
    'E35' = '  +5.38%  '
    'E36' = '  +0.19%  '
    'D37' = '0.0899'
    'E37' = '  +10.16%  '
    'D38' = '469.12'
    'E38' = '  +5.41%  '
    'D39' = '0.0428'
    'E39' = '  +3.94%  '
    'D40' = '3.04'
    'E40' = '  +7.93%  '
    'E41' = '  +4.16%  '
    'D42' = '3.070.66'
    'E42' = '  +1.31%  '
    'E43' = '  +0.57%  '
    'D44' = '2.44'
    'E44' = '  +7.55%  '
    'D45' = '0.284'
    'E45' = '  +3.70%  '
    'D46' = '28.86'
    'E46' = '  +3.50%  '
    'D47' = '0.0₃0591'
    'E47' = '  +13.99%  '
    'E48' = '  -0.04%  '
    'E49' = '  +0.91%  '
    'D50' = '2.26'
    'E50' = '  +6.37%  '
    'D51' = '120.01'
    'E51' = '  +1.58%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
